$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 122, shifting existing rows 122:130 down to 123:131
$ws.Rows(122).Insert()

# Populate the new row 122 with the new record's data
$ws.Cells.Item(122, 1).Value = 11
$ws.Cells.Item(122, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(122, 3).Value = "Bíobío"
$ws.Cells.Item(122, 4).Value = 44714
$ws.Cells.Item(122, 5).Value = 8
$ws.Cells.Item(122, 6).Value = 100112043
$ws.Cells.Item(122, 7).Value = "Pepino ensalada"
$ws.Cells.Item(122, 8).Value = "Sin especificar"
$ws.Cells.Item(122, 9).Value = "Primera"
$ws.Cells.Item(122, 10).Value = 150
$ws.Cells.Item(122, 11).Value = 19000
$ws.Cells.Item(122, 12).Value = 21000
$ws.Cells.Item(122, 13).Value = 20067
$ws.Cells.Item(122, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(122, 15).Value = "Región Metropolitana"
$ws.Cells.Item(122, 16).Value = 334
$ws.Cells.Item(122, 17).Value = 60
$ws.Cells.Item(122, 18).Value = "Hortaliza"
